$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.745.57'
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').Value = '2.473.82'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'318.84"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.23%  '
$ws.Range('D6').Value = "'92.94"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.02%  '
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').Value = "'0.0871"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +9.77%  '
$ws.Range('D11').Value = "'33.34"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.52%  '
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('D13').Value = '2.854.98'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('E14').Value = '  +0.65%  '
$ws.Range('D15').Value = "'15.66"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.09%  '
$ws.Range('D16').Value = '2.456.76'
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('D17').Value = "'0.796"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.30%  '
$ws.Range('D18').Value = '41.697.78'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').Value = "'6.46"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('D20').Value = '0.0₃0950'
$ws.Range('E20').Value = '  +0.68%  '
$ws.Range('D21').Value = "'71.06"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('D22').Value = "'11.31"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.67%  '
$ws.Range('D23').Value = "'240.56"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.45%  '
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('E28').Value = '  +0.81%  '
$ws.Range('E29').Value = '  +0.85%  '
$ws.Range('D30').Value = "'36.37"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.58%  '
$ws.Range('D31').Value = "'158.03"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.20%  '
$ws.Range('E32').Value = '  +0.72%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('E35').Value = '  +0.70%  '
$ws.Range('D36').Value = "'17.47"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.38%  '
$ws.Range('E37').Value = '  +4.27%  '
$ws.Range('E38').Value = '  +1.02%  '
$ws.Range('E39').Value = '  +1.96%  '
$ws.Range('E40').Value = '  +0.72%  '
$ws.Range('D41').Value = "'2.56"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.00%  '
$ws.Range('D42').Value = "'4.00"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').Value = '1.997.94'
$ws.Range('E43').Value = '  +2.65%  '
$ws.Range('D45').Value = "'18.90"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.41%  '
$ws.Range('E46').Value = '  +2.59%  '
$ws.Range('D47').Value = "'9.58"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.42%  '
$ws.Range('D48').Value = '2.711.77'
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('D49').Value = "'97.97"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.58%  '
$ws.Range('D50').Value = "'74.85"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.46%  '
$ws.Range('D51').Value = "'67.11"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.19%  '
